$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 (header): add P1=14, Q1=15, copying the style of the existing header cells ---
$ws.Range("O1").Copy() | Out-Null
$ws.Range("P1:Q1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Cells.Item(1, 16).Value = 14
$ws.Cells.Item(1, 17).Value = 15

# --- Row 2: replace all-zero values with the new pattern, and add P2/Q2 ---
$row2Values = @(2,2,2,1,1,1,2,2,2,1,2,2,2,1)
for ($i = 0; $i -lt $row2Values.Length; $i++) {
    $ws.Cells.Item(2, 2 + $i).Value = $row2Values[$i]
}
$ws.Cells.Item(2, 16).Value = 2
$ws.Cells.Item(2, 17).Value = 2

# --- Rows 3-25: flip I/K/M/O columns and append P/Q columns ---
for ($r = 3; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value = 2    # I -> 2
    $ws.Cells.Item($r, 11).Value = 1   # K -> 1
    $ws.Cells.Item($r, 13).Value = 2   # M -> 2
    $ws.Cells.Item($r, 15).Value = 1   # O -> 1
    $ws.Cells.Item($r, 16).Value = 2   # P = 2
    $ws.Cells.Item($r, 17).Value = 2   # Q = 2
}
